$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The quarterly database is laid out in columns E:N (10 trailing quarters).
# This update rolls the window forward by one quarter: the oldest quarter
# (column E, "Q2 ended 1399/06") drops off, every remaining quarter's data
# shifts one column to the left, and the newest quarter ("Q4 ended 1401/12")
# is appended in column N - both for the quarter-header labels and for every
# metric row below them.

$cols = @("E", "F", "G", "H", "I", "J", "K", "L", "M", "N")

# New quarter-header captions for columns E..N after the shift (row 8 and row 24).
$quarterLabels = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

$headerRows = @(8, 24)
foreach ($row in $headerRows) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $quarterLabels[$i]
    }
}

# Updated (re-priced) quarterly figures for columns E..N after the shift,
# one array per data row.
$rowValues = @{
    10 = @(207394, 315925, 84595, 455581, 553197, 501708, 435006, 406444, 787781, 702897)
    13 = @(705, -1283, 637, 50, 765, 1094, 0, 1211, -1211, 3254)
    15 = @(-139, 515, 76, 322, 144, 313, 355, 166, -68, 572)
    16 = @(92, 2354, 563, 3987, -4285, 9147, 5511, 7251, -12762, 23129)
    17 = @(8395, 27867, 16726, 23878, -14449, 47376, 15618, 38207, 17300, 20744)
    19 = @(-11575, 67242, 252896, -223817, 40562, 28237, 23808, 37312, 37376, 69262)
    20 = @(204872, 412620, 355493, 260001, 575934, 587875, 480298, 490591, 828416, 819858)
    26 = @(226, 234, 234, 211, 211, 234, 234, 234, 234, 469)
    27 = @(26, 24, 24, 20, 20, 24, 24, 24, 24, 44)
}

foreach ($row in $rowValues.Keys) {
    $values = $rowValues[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}
